$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F:V content of rows 15 and 16 ---------------------------------
$row15 = $ws.Range("F15:V15").Value()
$row16 = $ws.Range("F16:V16").Value()

$ws.Range("F15:V15").Value = $row16
$ws.Range("F16:V16").Value = $row15

# --- Append three new match rows (57, 58, 59) --------------------------------
$newRows = @(
    @(56, "morocco", "botola-pro", "2023-2024", 45234.66666666666, "Renaissance Zemamra", 1, "Moghreb Tetouan", 1, 2.12, "03/11/2023 04:13", 2.7, "04/11/2023 15:52", 2.9, "03/11/2023 04:13", 2.79, "04/11/2023 14:11", 3.48, "03/11/2023 04:13", 2.92, "04/11/2023 15:52", "https://www.betexplorer.com/football/morocco/botola-pro/renaissance-zemamra-moghreb-tetouan/dhLxi7MB/"),
    @(57, "morocco", "botola-pro", "2023-2024", 45234.76041666666, "Berkane", 0, "Hassania Agadir", 0, 1.35, "03/11/2023 12:42", 1.38, "04/11/2023 17:49", 4.19, "03/11/2023 12:42", 4.16, "04/11/2023 17:50", 8.09, "03/11/2023 12:42", 9.449999999999999, "04/11/2023 17:50", "https://www.betexplorer.com/football/morocco/botola-pro/berkane-hassania-agadir/AePtjRyI/"),
    @(58, "morocco", "botola-pro", "2023-2024", 45234.85416666666, "FAR Rabat", 4, "Youssoufia Berrechid", 1, 1.41, "03/11/2023 15:13", 1.26, "04/11/2023 20:20", 3.95, "03/11/2023 15:13", 4.99, "04/11/2023 20:25", 6.6, "03/11/2023 15:13", 13.48, "04/11/2023 20:25", "https://www.betexplorer.com/football/morocco/botola-pro/far-rabat-youssoufia-berrechid/2imAYb6p/")
)

$lastDataRow = 56
$r = 57
foreach ($row in $newRows) {
    # Clone the formatting (bold/bordered/centered index cell, datetime
    # number format on the match-date cell) from the last existing data
    # row before filling in this row's values.
    $ws.Range("A$lastDataRow`:E$lastDataRow").Copy()
    $ws.Range("A$r`:E$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $ws.Cells.Item($r, 15).Value = $row[14]
    $ws.Cells.Item($r, 16).Value = $row[15]
    $ws.Cells.Item($r, 17).Value = $row[16]
    $ws.Cells.Item($r, 18).Value = $row[17]
    $ws.Cells.Item($r, 19).Value = $row[18]
    $ws.Cells.Item($r, 20).Value = $row[19]
    $ws.Cells.Item($r, 21).Value = $row[20]
    $ws.Cells.Item($r, 22).Value = $row[21]

    $r = $r + 1
}
